$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 611
$ws.Range("I28").Value = 520.4
$ws.Range("J28").Value = 1970
$ws.Range("K28").Value = 520.4
$ws.Range("L28").Value = 1970
$ws.Range("M28").Value = -35.39999999999998
$ws.Range("N28").Value = -2940
# Row 33
$ws.Range("H33").Value = 364.13043
$ws.Range("I33").Value = 409.3684
$ws.Range("K33").Value = 409.3684
$ws.Range("M33").Value = -180.3684
# Row 40
$ws.Range("H40").Value = 402371.75
$ws.Range("I40").Value = 627268.4
$ws.Range("J40").Value = 2555.5557
$ws.Range("K40").Value = 627268.4
$ws.Range("L40").Value = 2555.5557
$ws.Range("M40").Value = -627093.4
$ws.Range("N40").Value = -2905.5557
# Row 43
$ws.Range("H43").Value = 2816.3333
$ws.Range("I43").Value = 2700
$ws.Range("J43").Value = 2874.5
$ws.Range("K43").Value = 2700
$ws.Range("L43").Value = 2874.5
$ws.Range("M43").Value = -2631
$ws.Range("N43").Value = -3012.5
# Row 64
$ws.Range("H64").Value = 6203.381
$ws.Range("I64").Value = 5893.625
$ws.Range("K64").Value = 5893.625
$ws.Range("M64").Value = -5645.625
# Row 67
$ws.Range("H67").Value = 6203.381
$ws.Range("I67").Value = 5893.625
$ws.Range("K67").Value = 5893.625
$ws.Range("M67").Value = -5035.625
# Row 125
$ws.Range("H125").Value = 5675.273
$ws.Range("I125").Value = 7633
$ws.Range("J125").Value = 4941.125
$ws.Range("K125").Value = 68697
$ws.Range("L125").Value = 44470.125
$ws.Range("M125").Value = -66237
$ws.Range("N125").Value = -49390.125
# Row 138
$ws.Range("H138").Value = 26320374
$ws.Range("J138").Value = 62509636
$ws.Range("L138").Value = 187528908
$ws.Range("N138").Value = -187539188

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5444.263
$ws.Range("I32").Value = 5610.943
$ws.Range("K32").Value = 5610.943
$ws.Range("M32").Value = -5323.943
# Row 61
$ws.Range("H61").Value = 2643.7932
$ws.Range("I61").Value = 2695.2273
$ws.Range("K61").Value = 2695.2273
$ws.Range("M61").Value = -2483.2273
# Row 63
$ws.Range("H63").Value = 3358.7144
$ws.Range("I63").Value = 3358.7144
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 3358.7144
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -2672.7144
$ws.Range("N63").ClearContents()
# Row 66
$ws.Range("H66").Value = 3358.7144
$ws.Range("I66").Value = 3358.7144
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 16793.572
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -13361.572
$ws.Range("N66").ClearContents()
# Row 88
$ws.Range("H88").Value = 9748.462
$ws.Range("I88").Value = 17292
$ws.Range("J88").Value = 3282.5715
$ws.Range("K88").Value = 17292
$ws.Range("L88").Value = 3282.5715
$ws.Range("M88").Value = -16886
$ws.Range("N88").Value = -4094.5715
# Row 91
$ws.Range("H91").Value = 9748.462
$ws.Range("I91").Value = 17292
$ws.Range("J91").Value = 3282.5715
$ws.Range("K91").Value = 17292
$ws.Range("L91").Value = 3282.5715
$ws.Range("M91").Value = -15888
$ws.Range("N91").Value = -6090.5715
# Row 122
$ws.Range("H122").Value = 2579.5454
$ws.Range("I122").Value = 1834.625
$ws.Range("J122").Value = 4566
$ws.Range("K122").Value = 5503.875
$ws.Range("L122").Value = 13698
$ws.Range("M122").Value = -3053.875
$ws.Range("N122").Value = -18598
# Row 136
$ws.Range("H136").Value = 2643.7932
$ws.Range("I136").Value = 2695.2273
$ws.Range("K136").Value = 8085.6819
$ws.Range("M136").Value = -5535.6819

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 61
$ws.Range("H61").Value = 59999
$ws.Range("J61").Value = 59999
$ws.Range("L61").Value = 59999
$ws.Range("N61").Value = -60625
# Row 134
$ws.Range("H134").Value = 5850.2
$ws.Range("I134").Value = 1106
$ws.Range("K134").Value = 3318
$ws.Range("M134").Value = -783

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2242
$ws.Range("I16").Value = 2052.5
$ws.Range("K16").Value = 2052.5
$ws.Range("M16").Value = -1765.5
# Row 33
$ws.Range("H33").Value = 805.1667
$ws.Range("I33").Value = 776.2
$ws.Range("J33").Value = 950
$ws.Range("K33").Value = 776.2
$ws.Range("L33").Value = 950
$ws.Range("M33").Value = -397.2
$ws.Range("N33").Value = -1708
# Row 113
$ws.Range("H113").Value = 2242
$ws.Range("I113").Value = 2052.5
$ws.Range("K113").Value = 2052.5
$ws.Range("M113").Value = 117.5
# Row 132
$ws.Range("H132").Value = 54349.31
$ws.Range("I132").Value = 59983.516
$ws.Range("K132").Value = 179950.548
$ws.Range("M132").Value = -177420.548

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 62
$ws.Range("H62").Value = 16668169
$ws.Range("J62").Value = 66666664
$ws.Range("L62").Value = 199999992
$ws.Range("N62").Value = -200001364
# Row 65
$ws.Range("H65").Value = 16668169
$ws.Range("J65").Value = 66666664
$ws.Range("L65").Value = 599999976
$ws.Range("N65").Value = -600006840
# Row 86
$ws.Range("H86").Value = 893.46155
$ws.Range("I86").Value = 594.2778
$ws.Range("J86").Value = 1566.625
$ws.Range("K86").Value = 1782.8334
$ws.Range("L86").Value = 4699.875
$ws.Range("M86").Value = -596.8334
$ws.Range("N86").Value = -7071.875
# Row 89
$ws.Range("H89").Value = 893.46155
$ws.Range("I89").Value = 594.2778
$ws.Range("J89").Value = 1566.625
$ws.Range("K89").Value = 5348.500199999999
$ws.Range("L89").Value = 14099.625
$ws.Range("M89").Value = 579.4998000000005
$ws.Range("N89").Value = -25955.625
# Row 102
$ws.Range("H102").Value = 2949.2
$ws.Range("J102").Value = 4999
$ws.Range("L102").Value = 14997
$ws.Range("N102").Value = -19865
# Row 109
$ws.Range("H109").Value = 9604
$ws.Range("J109").Value = 1998
$ws.Range("L109").Value = 5994
$ws.Range("N109").Value = -8074
# Row 121
$ws.Range("H121").Value = 348.57144
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 93
$ws.Range("H93").Value = 33976
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 33976
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 33976
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -37720
# Row 94
$ws.Range("H94").Value = 28833
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 28833
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 28833
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -30185
# Row 126
$ws.Range("H126").Value = 4247.375
$ws.Range("I126").Value = 3800
$ws.Range("J126").Value = 4993
$ws.Range("K126").Value = 11400
$ws.Range("L126").Value = 14979
$ws.Range("M126").Value = -8930
$ws.Range("N126").Value = -19919

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1585.5294
$ws.Range("I22").Value = 1250.3846
$ws.Range("J22").Value = 2674.75
$ws.Range("K22").Value = 1250.3846
$ws.Range("L22").Value = 2674.75
$ws.Range("M22").Value = -955.3846000000001
$ws.Range("N22").Value = -3264.75
# Row 27
$ws.Range("H27").Value = 1585.5294
$ws.Range("I27").Value = 1250.3846
$ws.Range("J27").Value = 2674.75
$ws.Range("K27").Value = 1250.3846
$ws.Range("L27").Value = 2674.75
$ws.Range("M27").Value = -1143.3846
$ws.Range("N27").Value = -2888.75
# Row 40
$ws.Range("H40").Value = 1913.8334
$ws.Range("I40").Value = 1633.2727
$ws.Range("K40").Value = 1633.2727
$ws.Range("M40").Value = -1497.2727
# Row 46
$ws.Range("H46").Value = 2349.9
$ws.Range("I46").Value = 1000
$ws.Range("K46").Value = 1000
$ws.Range("M46").Value = -812
# Row 69
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
# Row 72
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
# Row 82
$ws.Range("H82").Value = 3108.7727
$ws.Range("I82").Value = 2606.1333
$ws.Range("J82").Value = 4185.857
$ws.Range("K82").Value = 2606.1333
$ws.Range("L82").Value = 4185.857
$ws.Range("M82").Value = -2245.1333
$ws.Range("N82").Value = -4907.857
# Row 85
$ws.Range("H85").Value = 3108.7727
$ws.Range("I85").Value = 2606.1333
$ws.Range("J85").Value = 4185.857
$ws.Range("K85").Value = 2606.1333
$ws.Range("L85").Value = 4185.857
$ws.Range("M85").Value = -1358.1333
$ws.Range("N85").Value = -6681.857
# Row 122
$ws.Range("H122").Value = 4317.125
$ws.Range("I122").Value = 4077.0833
$ws.Range("J122").Value = 5037.25
$ws.Range("K122").Value = 12231.2499
$ws.Range("L122").Value = 15111.75
$ws.Range("M122").Value = -9781.249899999999
$ws.Range("N122").Value = -20011.75
# Row 132
$ws.Range("H132").Value = 3857.1428
# Row 136
$ws.Range("H136").Value = 10047.5
$ws.Range("I136").Value = 5710.7144
$ws.Range("K136").Value = 17132.1432
$ws.Range("M136").Value = -14582.1432

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 5502.476
$ws.Range("J136").Value = 7510.4443
$ws.Range("L136").Value = 22531.3329
$ws.Range("N136").Value = -27631.3329
